$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current "room" column (D) so that
# new columns D,E hold row_dist/seat_dist, and room/seat_last shift to F,G.
$ws.Range("D1:E1").EntireColumn.Insert()

$ws.Range("D1").Value = "row_dist"
$ws.Range("E1").Value = "seat_dist"

$ws.Range("E2").Value = "0,5"
$ws.Range("D2").Value = "1,2"

$ws.Range("D1:E1").Font.Bold = $true

$ws.Range("D1:E2").NumberFormat = "@"
$ws.Columns("D:E").ColumnWidth = $ws.Columns("C:C").ColumnWidth

$ws.Range("D3").Select()
